# Auto-generated edit script: updates Leve profit-calculation columns (H-N)
# across multiple profession sheets to reflect refreshed market-price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 48376.75
$ws.Range("I62").Value = 64812.41
$ws.Range("K62").Value = 64812.41
$ws.Range("M62").Value = -64188.41
$ws.Range("H65").Value = 48376.75
$ws.Range("I65").Value = 64812.41
$ws.Range("K65").Value = 324062.05
$ws.Range("M65").Value = -320942.05
$ws.Range("H86").Value = 3779.3076
$ws.Range("I86").Value = 3234.2354
$ws.Range("J86").Value = 4808.8887
$ws.Range("K86").Value = 3234.2354
$ws.Range("L86").Value = 4808.8887
$ws.Range("M86").Value = -2111.2354
$ws.Range("N86").Value = -7054.8887
$ws.Range("H89").Value = 3779.3076
$ws.Range("I89").Value = 3234.2354
$ws.Range("J89").Value = 4808.8887
$ws.Range("K89").Value = 16171.177
$ws.Range("L89").Value = 24044.4435
$ws.Range("M89").Value = -10555.177
$ws.Range("N89").Value = -35276.4435
$ws.Range("H107").Value = 1669.5555
$ws.Range("I107").Value = 647.36365
$ws.Range("J107").Value = 3275.8572
$ws.Range("K107").Value = 647.36365
$ws.Range("L107").Value = 3275.8572
$ws.Range("M107").Value = 1272.63635
$ws.Range("N107").Value = -7115.8572
$ws.Range("H127").Value = 5000
$ws.Range("I127").Value = 5000
$ws.Range("K127").Value = 15000
$ws.Range("M127").Value = -10040
$ws.Range("H132").Value = 13932776
$ws.Range("I132").Value = 15199283
$ws.Range("K132").Value = 45597849
$ws.Range("M132").Value = -45595319
$ws.Range("H135").Value = 3662.75
$ws.Range("I135").Value = 3677.5454
$ws.Range("J135").Value = 3500
$ws.Range("K135").Value = 33097.9086
$ws.Range("L135").Value = 31500
$ws.Range("M135").Value = -30562.9086
$ws.Range("N135").Value = -36570
$ws.Range("H138").Value = 4197.6313
$ws.Range("J138").Value = 4290.0576
$ws.Range("L138").Value = 12870.1728
$ws.Range("N138").Value = -23150.1728

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1736.7368
$ws.Range("I32").Value = 1511.9718
$ws.Range("K32").Value = 1511.9718
$ws.Range("M32").Value = -1224.9718
$ws.Range("H61").Value = 3226.879
$ws.Range("I61").Value = 2878.862
$ws.Range("J61").Value = 5750
$ws.Range("K61").Value = 2878.862
$ws.Range("L61").Value = 5750
$ws.Range("M61").Value = -2666.862
$ws.Range("N61").Value = -6174
$ws.Range("H74").Value = 51205.953
$ws.Range("I74").Value = 51076.6
$ws.Range("J74").Value = 52499.5
$ws.Range("K74").Value = 51076.6
$ws.Range("L74").Value = 52499.5
$ws.Range("M74").Value = -50202.6
$ws.Range("N74").Value = -54247.5
$ws.Range("H77").Value = 51205.953
$ws.Range("I77").Value = 51076.6
$ws.Range("J77").Value = 52499.5
$ws.Range("K77").Value = 255383
$ws.Range("L77").Value = 262497.5
$ws.Range("M77").Value = -251015
$ws.Range("N77").Value = -271233.5
$ws.Range("H110").Value = 13194.15
$ws.Range("I110").Value = 21709.445
$ws.Range("J110").Value = 6227.091
$ws.Range("K110").Value = 21709.445
$ws.Range("L110").Value = 6227.091
$ws.Range("M110").Value = -19664.445
$ws.Range("N110").Value = -10317.091
$ws.Range("H132").Value = 6493.077
$ws.Range("I132").Value = 6858.857
$ws.Range("K132").Value = 20576.571
$ws.Range("M132").Value = -18046.571
$ws.Range("H136").Value = 3226.879
$ws.Range("I136").Value = 2878.862
$ws.Range("J136").Value = 5750
$ws.Range("K136").Value = 8636.585999999999
$ws.Range("L136").Value = 17250
$ws.Range("M136").Value = -6086.585999999999
$ws.Range("N136").Value = -22350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3640.5925
$ws.Range("I99").Value = 1732.238
$ws.Range("J99").Value = 10319.833
$ws.Range("K99").Value = 1732.238
$ws.Range("L99").Value = 10319.833
$ws.Range("M99").Value = -234.2380000000001
$ws.Range("N99").Value = -13315.833
$ws.Range("H107").Value = 1706.2727
$ws.Range("I107").Value = 1415.4736
$ws.Range("J107").Value = 2100.9285
$ws.Range("K107").Value = 1415.4736
$ws.Range("L107").Value = 2100.9285
$ws.Range("M107").Value = 504.5264
$ws.Range("N107").Value = -5940.9285
$ws.Range("H134").Value = 4689.8
$ws.Range("I134").Value = 4225
$ws.Range("J134").Value = 4999.6665
$ws.Range("K134").Value = 12675
$ws.Range("L134").Value = 14998.9995
$ws.Range("M134").Value = -10140
$ws.Range("N134").Value = -20068.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3961133.5
$ws.Range("I31").Value = 5276858
$ws.Range("J31").Value = 13960.25
$ws.Range("K31").Value = 5276858
$ws.Range("L31").Value = 13960.25
$ws.Range("M31").Value = -5276563
$ws.Range("N31").Value = -14550.25
$ws.Range("H34").Value = 3961133.5
$ws.Range("I34").Value = 5276858
$ws.Range("J34").Value = 13960.25
$ws.Range("K34").Value = 5276858
$ws.Range("L34").Value = 13960.25
$ws.Range("M34").Value = -5276656
$ws.Range("N34").Value = -14364.25
$ws.Range("H107").Value = 689.6786
$ws.Range("I107").Value = 652.4400000000001
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 652.4400000000001
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1267.56
$ws.Range("N107").Value = -4840
$ws.Range("H132").Value = 4466.375
$ws.Range("I132").Value = 4507.7827
$ws.Range("J132").Value = 3514
$ws.Range("K132").Value = 13523.3481
$ws.Range("L132").Value = 10542
$ws.Range("M132").Value = -10993.3481
$ws.Range("N132").Value = -15602
$ws.Range("H134").Value = 27549.87
$ws.Range("I134").Value = 27549.87
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 82649.61
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -80114.61
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1217
$ws.Range("N32").ClearContents()
$ws.Range("H46").Value = 546.75
$ws.Range("I46").Value = 274.8
$ws.Range("K46").Value = 824.4000000000001
$ws.Range("M46").Value = -733.4000000000001
$ws.Range("H68").Value = 2385301.2
$ws.Range("I68").Value = 4440
$ws.Range("J68").Value = 2635918
$ws.Range("K68").Value = 13320
$ws.Range("L68").Value = 7907754
$ws.Range("M68").Value = -12509
$ws.Range("N68").Value = -7909376
$ws.Range("H71").Value = 2385301.2
$ws.Range("I71").Value = 4440
$ws.Range("J71").Value = 2635918
$ws.Range("K71").Value = 39960
$ws.Range("L71").Value = 23723262
$ws.Range("M71").Value = -35904
$ws.Range("N71").Value = -23731374
$ws.Range("H107").Value = 22224862
$ws.Range("I107").Value = 111111630
$ws.Range("J107").Value = 3169.4167
$ws.Range("K107").Value = 333334890
$ws.Range("L107").Value = 9508.250100000001
$ws.Range("M107").Value = -333332970
$ws.Range("N107").Value = -13348.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2826.2703
$ws.Range("I132").Value = 2837.5881
$ws.Range("J132").Value = 2698
$ws.Range("K132").Value = 8512.764299999999
$ws.Range("L132").Value = 8094
$ws.Range("M132").Value = -5982.764299999999
$ws.Range("N132").Value = -13154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3707.2727
$ws.Range("I132").Value = 3325.7144
$ws.Range("J132").Value = 4375
$ws.Range("K132").Value = 9977.143199999999
$ws.Range("L132").Value = 13125
$ws.Range("M132").Value = -7447.143199999999
$ws.Range("N132").Value = -18185
$ws.Range("H136").Value = 3549.6155
$ws.Range("I136").Value = 3549.6155
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10648.8465
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8098.8465
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 142893470
$ws.Range("J125").Value = 142893470
$ws.Range("L125").Value = 142893470
$ws.Range("N125").Value = -142903310
